# "wrap up referenced parsing"
# Rename the ".ErrorCasesParsing" sheet to "ErrorCasesParsing" (drop the
# leading dot now that the sheet's parsing work is wired up / referenced
# elsewhere), and leave the sheet's selection where the author last left
# it (cell D13) instead of the old scratch selection (U46).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(".ErrorCasesParsing")
$ws.Name = "ErrorCasesParsing"

$ws.Activate()
$ws.Range("D13").Select()
